$d = $word.ActiveDocument

$replacements = @(
    @{old = "232÷2=116, 0"; new = "516÷8=64, 4"},
    @{old = "435÷3=145, 0"; new = "808÷5=161, 3"},
    @{old = "920÷9=102, 2"; new = "994÷3=331, 1"},
    @{old = "206÷4=51, 2"; new = "865÷8=108, 1"},
    @{old = "987÷7=141, 0"; new = "273÷6=45, 3"},
    @{old = "149÷5=29, 4"; new = "368÷8=46, 0"},
    @{old = "637÷6=106, 1"; new = "400÷9=44, 4"},
    @{old = "125÷8=15, 5"; new = "196÷5=39, 1"},
    @{old = "857÷5=171, 2"; new = "755÷4=188, 3"},
    @{old = "461÷9=51, 2"; new = "531÷3=177, 0"},
    @{old = "665÷6=110, 5"; new = "115÷4=28, 3"},
    @{old = "251÷3=83, 2"; new = "897÷2=448, 1"},
    @{old = "436÷4=109, 0"; new = "366÷7=52, 2"},
    @{old = "232÷4=58, 0"; new = "696÷9=77, 3"},
    @{old = "398÷6=66, 2"; new = "284÷4=71, 0"},
    @{old = "323÷3=107, 2"; new = "897÷9=99, 6"},
    @{old = "887÷9=98, 5"; new = "105÷5=21, 0"},
    @{old = "655÷9=72, 7"; new = "677÷2=338, 1"},
    @{old = "122÷3=40, 2"; new = "860÷2=430, 0"},
    @{old = "956÷8=119, 4"; new = "368÷3=122, 2"},
    @{old = "367÷6=61, 1"; new = "305÷6=50, 5"},
    @{old = "456÷9=50, 6"; new = "445÷3=148, 1"},
    @{old = "380÷5=76, 0"; new = "837÷5=167, 2"},
    @{old = "963÷4=240, 3"; new = "742÷6=123, 4"},
    @{old = "671÷5=134, 1"; new = "194÷6=32, 2"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
